# Add a "Live DEMO Link" textbox with a hyperlink to slide 4
# (sldId="275", creationId "{4288587416}" / cId 4288587416).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Create the new textbox (PowerPoint will name/number it the same way the
# real app would: "TextBox 6" / shape id 7, since 6 shapes already exist
# after the insert).
$tb = $s.Shapes.AddTextbox(1, 100, 100, 200, 50)

# --- position & size -------------------------------------------------
# Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU)
# and are rounded to a 32-bit float internally before being converted
# back to EMU, so values are nudged by a hair above the exact quotient
# to avoid landing 1 EMU low after truncation.
$tb.Left = 106.90834645669291
$tb.Top = 500.399857519685
$tb.Width = 481.46645669291337
$tb.Height = 29.081259842519685

# --- shape fill / text box behaviour ----------------------------------
$tb.Fill.Visible = $false

$tf = $tb.TextFrame
$tf.WordWrap = $false
$tf.AutoSize = 1          # ppAutoSizeShapeToFitText -> <a:spAutoFit/>

$tr = $tf.TextRange
$fullText = "Live DEMO Link: https://youtu.be/OkruHlKifrg?feature=shared "
$tr.Text = $fullText

$tr.ParagraphFormat.Alignment = 2   # ppAlignCenter -> <a:pPr algn="ctr"/>

# Run 1: "Live DEMO Link" (bold)
$run1 = $tr.Characters(1, 14)
$run1.Font.Bold = $true

# Run 2: ": "
$run2 = $tr.Characters(15, 2)

# Run 3: the URL, turned into a hyperlink
$run3 = $tr.Characters(17, 43)
$run3.ActionSettings.Item(1).Hyperlink.Address = "https://youtu.be/OkruHlKifrg?feature=shared"

# Run 4: trailing space
$run4 = $tr.Characters(60, 1)
